# Insert a new data row at row 19 (shifts existing rows 19-82 down to 20-83)
# and populate it with the new record described by the commit diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 19; this pushes the old row 19
# (and everything below it) down by one, matching the diff's row shift.
$ws.Rows.Item(19).Insert()

# Populate the newly inserted row 19 with the new record's values.
$ws.Cells.Item(19, 1).Value2 = 10
$ws.Cells.Item(19, 2).Value2 = "Vega Modelo de Temuco"
$ws.Cells.Item(19, 3).Value2 = "La Araucanía"
$ws.Cells.Item(19, 4).Value2 = 44910
$ws.Cells.Item(19, 5).Value2 = 9
$ws.Cells.Item(19, 6).Value2 = 100112030
$ws.Cells.Item(19, 7).Value2 = "Poroto granado"
$ws.Cells.Item(19, 8).Value2 = "Sin especificar"
$ws.Cells.Item(19, 9).Value2 = "Primera"
$ws.Cells.Item(19, 10).Value2 = 15
$ws.Cells.Item(19, 11).Value2 = 50000
$ws.Cells.Item(19, 12).Value2 = 50000
$ws.Cells.Item(19, 13).Value2 = 50000
$ws.Cells.Item(19, 14).Value2 = "$/saco 25 kilos"
$ws.Cells.Item(19, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(19, 16).Value2 = 2000
$ws.Cells.Item(19, 17).Value2 = 25
$ws.Cells.Item(19, 18).Value2 = "Hortaliza"
